$d = $word.ActiveDocument

# 1. Fix "mailformed" -> "mail formed" in the XML paragraph
$d.Content.Find.Execute("mailformed", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "mail formed", 2) | Out-Null

# 2. Move the Word-managed "_GoBack" bookmark so it spans from the start of the
#    "XML" heading paragraph through the end of the paragraph that was just edited
#    (mirrors Word's behaviour of tracking the span of the most recent edits).
$startPara = $d.Paragraphs.Item(22)
$endPara = $d.Paragraphs.Item(25)
$goBackRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
